# Auto-generated Excel COM-interop edit script
# Applies updated market-price-derived figures (columns H-N) across all 8 sheets
# per the scheduled runner's refreshed pricing pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 8783.799999999999
$ws.Range("I29").Value = 1614.8334
$ws.Range("K29").Value = 4844.5002
$ws.Range("M29").Value = -4563.5002
$ws.Range("H51").Value = 6132.7896
$ws.Range("I51").Value = 4663.8887
$ws.Range("J51").Value = 7454.8
$ws.Range("K51").Value = 4663.8887
$ws.Range("L51").Value = 7454.8
$ws.Range("M51").Value = -4179.8887
$ws.Range("N51").Value = -8422.799999999999
$ws.Range("H58").Value = 65.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H129").Value = 764.9167
$ws.Range("I129").Value = 764.9167
$ws.Range("K129").Value = 2294.7501
$ws.Range("M129").Value = 2705.2499
$ws.Range("H135").Value = 2418.7
$ws.Range("I135").Value = 2585
$ws.Range("K135").Value = 23265
$ws.Range("M135").Value = -20730
$ws.Range("H137").Value = 2834.5862
$ws.Range("I137").Value = 1523.1666
$ws.Range("J137").Value = 3176.6956
$ws.Range("K137").Value = 4569.4998
$ws.Range("L137").Value = 9530.086800000001
$ws.Range("M137").Value = -2019.4998
$ws.Range("N137").Value = -14630.0868
$ws.Range("H138").Value = 4778.347
$ws.Range("I138").Value = 2362.85
$ws.Range("J138").Value = 5707.385
$ws.Range("K138").Value = 7088.549999999999
$ws.Range("L138").Value = 17122.155
$ws.Range("M138").Value = -1948.549999999999
$ws.Range("N138").Value = -27402.155
$ws.Range("H140").Value = 61266.668
$ws.Range("J140").Value = 60175
$ws.Range("L140").Value = 60175
$ws.Range("N140").Value = -70535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2355
$ws.Range("I32").Value = 1932.1017
$ws.Range("K32").Value = 1932.1017
$ws.Range("M32").Value = -1645.1017
$ws.Range("H61").Value = 15114.6
$ws.Range("I61").Value = 16146.866
$ws.Range("J61").Value = 14495.24
$ws.Range("K61").Value = 16146.866
$ws.Range("L61").Value = 14495.24
$ws.Range("M61").Value = -15934.866
$ws.Range("N61").Value = -14919.24
$ws.Range("H74").Value = 3521.1538
$ws.Range("I74").Value = 2070.6155
$ws.Range("J74").Value = 4246.423
$ws.Range("K74").Value = 2070.6155
$ws.Range("L74").Value = 4246.423
$ws.Range("M74").Value = -1196.6155
$ws.Range("N74").Value = -5994.423
$ws.Range("H77").Value = 3521.1538
$ws.Range("I77").Value = 2070.6155
$ws.Range("J77").Value = 4246.423
$ws.Range("K77").Value = 10353.0775
$ws.Range("L77").Value = 21232.115
$ws.Range("M77").Value = -5985.077499999999
$ws.Range("N77").Value = -29968.115
$ws.Range("H102").Value = 429300.2
$ws.Range("I102").Value = 443123.44
$ws.Range("K102").Value = 443123.44
$ws.Range("M102").Value = -441501.44
$ws.Range("H122").Value = 837167.6
$ws.Range("I122").Value = 1113662.4
$ws.Range("J122").Value = 7683.3335
$ws.Range("K122").Value = 3340987.2
$ws.Range("L122").Value = 23050.0005
$ws.Range("M122").Value = -3338537.2
$ws.Range("N122").Value = -27950.0005
$ws.Range("H136").Value = 15114.6
$ws.Range("I136").Value = 16146.866
$ws.Range("J136").Value = 14495.24
$ws.Range("K136").Value = 48440.598
$ws.Range("L136").Value = 43485.72
$ws.Range("M136").Value = -45890.598
$ws.Range("N136").Value = -48585.72

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 29250
$ws.Range("J9").Value = 29250
$ws.Range("L9").Value = 29250
$ws.Range("N9").Value = -29586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 5000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -4887
$ws.Range("N6").ClearContents()
$ws.Range("H9").Value = 46507
$ws.Range("J9").Value = 46507
$ws.Range("L9").Value = 46507
$ws.Range("N9").Value = -46843
$ws.Range("H31").Value = 29415946
$ws.Range("I31").Value = 90911144
$ws.Range("J31").Value = 5201.391
$ws.Range("K31").Value = 90911144
$ws.Range("L31").Value = 5201.391
$ws.Range("M31").Value = -90910849
$ws.Range("N31").Value = -5791.391
$ws.Range("H34").Value = 29415946
$ws.Range("I34").Value = 90911144
$ws.Range("J34").Value = 5201.391
$ws.Range("K34").Value = 90911144
$ws.Range("L34").Value = 5201.391
$ws.Range("M34").Value = -90910942
$ws.Range("N34").Value = -5605.391
$ws.Range("H58").Value = 590004.2
$ws.Range("I58").Value = 834599.2
$ws.Range("J58").Value = 2976.3
$ws.Range("K58").Value = 834599.2
$ws.Range("L58").Value = 2976.3
$ws.Range("M58").Value = -834396.2
$ws.Range("N58").Value = -3382.3
$ws.Range("H99").Value = 17485.5
$ws.Range("I99").Value = 57506
$ws.Range("J99").Value = 9481.4
$ws.Range("K99").Value = 57506
$ws.Range("L99").Value = 9481.4
$ws.Range("M99").Value = -56008
$ws.Range("N99").Value = -12477.4
$ws.Range("H126").Value = 17485.5
$ws.Range("I126").Value = 57506
$ws.Range("J126").Value = 9481.4
$ws.Range("K126").Value = 172518
$ws.Range("L126").Value = 28444.2
$ws.Range("M126").Value = -170048
$ws.Range("N126").Value = -33384.2
$ws.Range("H136").Value = 590004.2
$ws.Range("I136").Value = 834599.2
$ws.Range("J136").Value = 2976.3
$ws.Range("K136").Value = 2503797.6
$ws.Range("L136").Value = 8928.900000000001
$ws.Range("M136").Value = -2501247.6
$ws.Range("N136").Value = -14028.9
$ws.Range("H141").Value = 150201
$ws.Range("J141").Value = 150201
$ws.Range("L141").Value = 150201
$ws.Range("N141").Value = -160561

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 226
$ws.Range("I14").Value = 226
$ws.Range("K14").Value = 678
$ws.Range("M14").Value = -505
$ws.Range("H107").Value = 497751
$ws.Range("I107").Value = 1445.25
$ws.Range("J107").Value = 795534.4399999999
$ws.Range("K107").Value = 4335.75
$ws.Range("L107").Value = 2386603.32
$ws.Range("M107").Value = -2415.75
$ws.Range("N107").Value = -2390443.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 310121.6
$ws.Range("J122").Value = 6457.2144
$ws.Range("L122").Value = 19371.6432
$ws.Range("N122").Value = -24271.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3519.41
$ws.Range("I7").Value = 3351.1948
$ws.Range("J7").Value = 4082.5652
$ws.Range("K7").Value = 3351.1948
$ws.Range("L7").Value = 4082.5652
$ws.Range("M7").Value = -3239.1948
$ws.Range("N7").Value = -4306.5652
$ws.Range("H40").Value = 12349915
$ws.Range("I40").Value = 13892654
$ws.Range("K40").Value = 13892654
$ws.Range("M40").Value = -13892518
$ws.Range("H122").Value = 5749.5884
$ws.Range("I122").Value = 3962.8572
$ws.Range("K122").Value = 11888.5716
$ws.Range("M122").Value = -9438.571599999999
$ws.Range("H126").Value = 3519.41
$ws.Range("I126").Value = 3351.1948
$ws.Range("J126").Value = 4082.5652
$ws.Range("K126").Value = 10053.5844
$ws.Range("L126").Value = 12247.6956
$ws.Range("M126").Value = -7583.5844
$ws.Range("N126").Value = -17187.6956
$ws.Range("H133").Value = 119988.5
$ws.Range("J133").Value = 119988.5
$ws.Range("L133").Value = 119988.5
$ws.Range("N133").Value = -125048.5
$ws.Range("H136").Value = 3831.7576
$ws.Range("I136").Value = 3830.3608
$ws.Range("J136").Value = 3899.5
$ws.Range("K136").Value = 11491.0824
$ws.Range("L136").Value = 11698.5
$ws.Range("M136").Value = -8941.082399999999
$ws.Range("N136").Value = -16798.5
$ws.Range("H140").Value = 122803.86
$ws.Range("J140").Value = 122803.86
$ws.Range("L140").Value = 122803.86
$ws.Range("N140").Value = -133163.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4786996
$ws.Range("I62").Value = 8774339
$ws.Range("J62").Value = 2183.8
$ws.Range("K62").Value = 8774339
$ws.Range("L62").Value = 2183.8
$ws.Range("M62").Value = -8773715
$ws.Range("N62").Value = -3431.8
$ws.Range("H65").Value = 4786996
$ws.Range("I65").Value = 8774339
$ws.Range("J65").Value = 2183.8
$ws.Range("K65").Value = 43871695
$ws.Range("L65").Value = 10919
$ws.Range("M65").Value = -43868575
$ws.Range("N65").Value = -17159
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H126").Value = 2977
$ws.Range("I126").Value = 2498.6
$ws.Range("K126").Value = 7495.799999999999
$ws.Range("M126").Value = -5025.799999999999
$ws.Range("H132").Value = 18529564
$ws.Range("I132").Value = 2537314
$ws.Range("K132").Value = 7611942
$ws.Range("M132").Value = -7609412
$ws.Range("H136").Value = 9133.968999999999
$ws.Range("I136").Value = 3436.1667
$ws.Range("J136").Value = 9938.364
$ws.Range("K136").Value = 10308.5001
$ws.Range("L136").Value = 29815.092
$ws.Range("M136").Value = -7758.500100000001
$ws.Range("N136").Value = -34915.092
